$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Former row 2 ("Adaboost") is removed; subsequent rows shift up and pick up
# refreshed metric values plus updated "Model Details" strings.
$ws.Rows("2:2").Delete()

# Refresh "Model Details" (column B) text for the remaining models
$ws.Range("B2").Value = "MLPClassifier(batch_size=32, early_stopping=True, max_iter=512)"
$ws.Range("B3").Value = "MultinomialNB()"
$ws.Range("B4").Value = "RandomForestClassifier(max_features='sqrt', n_jobs=8)"
$ws.Range("B5").Value = "SVC(degree=1, max_iter=1024)"

# Refresh Mean Acc (C) and Mean F1 (E) values; Std Acc (D) / Std F1 (F) stay 0
$ws.Range("C2").Value = 0.7885835095137421
$ws.Range("E2").Value = 0.7075336367233874

$ws.Range("C3").Value = 0.5940803382663847
$ws.Range("E3").Value = 0.5777571136321369

$ws.Range("C4").Value = 0.7822410147991543
$ws.Range("E4").Value = 0.6978273139447617

$ws.Range("C5").Value = 0.7864693446088795
$ws.Range("E5").Value = 0.6999428438632523

$ws.Range("C6").Value = 0.7463002114164905
$ws.Range("E6").Value = 0.427360774818402
